$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# NOTE: straight double-quote characters placed inside the *replacement* text
# get auto "smart-quoted" by this runtime's Find/Replace (mirroring Word's
# AutoFormat behaviour). To avoid mangling pre-existing straight quotes that
# must stay untouched, every search/replace pair below is scoped so that it
# never spans an existing quote character; new quote marks are always
# inserted using the literal curly characters "\u201c"/"\u201d" already.

# 1. Descripcion row: `" en la vista VerVentasView` -> `" en la vista “VerVentasView”`
Replace-Text 'vista VerVentasView' 'vista “VerVentasView”'

# 2. Disparador row: `” en VerVentasView` -> `” en “VerVentasView”`
Replace-Text 'Cancelar Venta” en VerVentasView' 'Cancelar Venta” en “VerVentasView”'

# 3. Precondiciones row PRE-04: add quotes around VerVentasView
Replace-Text 'una venta en la ventana VerVentasView.' 'una venta en la ventana “VerVentasView”.'

# 4. Flujo normal step 1: `muestra la venta ConfirmacionView` -> `muestra la venta “ConfirmacionView”`
Replace-Text 'muestra la venta ConfirmacionView' 'muestra la venta “ConfirmacionView”'

# 5. Flujo normal step 3 (EX-01 branch / InfoView paragraph) - full rewrite
Replace-Text 'El sistema muestra una venta InfoView con el mensaje ' 'El sistema cierra la ventana “ConfirmacionView” y muestra la ventana “InfoView” con el mensaje '

# 6. Flujo normal step 4: closing windows + showing VerVentasView
Replace-Text 'El sistema cierra la ventana InfoView, ConfirmacionView y muestra la ventana VerVentasView.' 'El sistema cierra la ventana “InfoView” y muestra la ventana “VerVentasView”.'

# 7. Flujo alterno: `El sistema cierra la ventana ConfirmacionView.`
Replace-Text 'El sistema cierra la ventana ConfirmacionView.' 'El sistema cierra la ventana “ConfirmacionView”.'

# 8. Excepciones: `la ventana ErrorView con el mensaje`
Replace-Text 'la ventana ErrorView con el mensaje' 'la ventana “ErrorView” con el mensaje'

# 9. Excepciones: `El sistema cierra las ventanas ErrorView, ConfirmacionView y VerVentasView.`
Replace-Text 'El sistema cierra las ventanas ErrorView, ConfirmacionView y VerVentasView.' 'El sistema cierra las ventanas “ErrorView” y “VerVentasView”.'
